# Edit: "11 May - Noche"
# Swap the Materia (E) / Docente (F) assignments for several rows on the
# "Blancos" sheet so each maps to a different (but still valid, already
# paired) subject/teacher combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

$ws.Range('E2').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F2').Value = 'Hernández Mendoza Delfina'
$ws.Range('E3').Value = 'TEMAS DE FÍSICA'
$ws.Range('F3').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E4').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F4').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E5').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F5').Value = 'Ortega Valle Manuel'
$ws.Range('E6').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F6').Value = 'Velasco Sanchez David'
$ws.Range('E7').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F7').Value = 'Velasco Sanchez David'
$ws.Range('E9').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F9').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E12').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F12').Value = 'Velasco Sanchez David'
$ws.Range('E13').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F13').Value = 'Hernández Mendoza Delfina'
$ws.Range('E17').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F17').Value = 'Velasco Sanchez David'
$ws.Range('E18').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F18').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E20').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F20').Value = 'Hernández Mendoza Delfina'
$ws.Range('E21').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F21').Value = 'Velasco Sanchez David'
$ws.Range('E23').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F23').Value = 'Hernández Mendoza Delfina'
$ws.Range('E24').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F24').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E32').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F32').Value = 'Velasco Sanchez David'
$ws.Range('E33').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F33').Value = 'Ortega Valle Manuel'
$ws.Range('E34').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F34').Value = 'Hernández Mendoza Delfina'
$ws.Range('E35').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F35').Value = 'Ortega Valle Manuel'
$ws.Range('E36').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F36').Value = 'Velasco Sanchez David'
$ws.Range('E39').Value = 'TEMAS DE FÍSICA'
$ws.Range('F39').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E40').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F40').Value = 'Hernández Mendoza Delfina'
$ws.Range('E41').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F41').Value = 'Ortega Valle Manuel'
$ws.Range('E42').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F42').Value = 'Ortega Valle Manuel'
$ws.Range('E43').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F43').Value = 'Hernández Mendoza Delfina'
$ws.Range('E44').Value = 'REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA'
$ws.Range('F44').Value = 'Zarate Amezcua Eladio Jorge'
$ws.Range('E45').Value = 'PROBABILIDAD Y ESTADÍSTICA'
$ws.Range('F45').Value = 'Velasco Sanchez David'
$ws.Range('E46').Value = 'TEMAS DE FÍSICA'
$ws.Range('F46').Value = 'Duran Amezcua Maria Angelica'
$ws.Range('E56').Value = 'TEMAS DE FILOSOFÍA'
$ws.Range('F56').Value = 'Hernández Mendoza Delfina'
$ws.Range('E57').Value = 'MATEMÁTICAS APLICADAS'
$ws.Range('F57').Value = 'Ortega Valle Manuel'
